$d = $word.ActiveDocument

$d.Content.Find.Execute("LEONARDO SILVERIO FERREIRA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MANOEL JEFETE DA SILVA TENORIO", 2)
